# Update countries & provincias Spain
# - Refresh the daily COVID-19 counters for a set of countries.
# - Two countries change rank/order after the refresh, which swaps the
#   labels shown on a couple of rows: Chile/Japon and
#   Madagascar/Liberia/Etiopia.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Country-label swaps caused by the re-sort of the refreshed data ---
$ws.Range("A31").Value = "Chile"
$ws.Range("A32").Value = "Japon"
$ws.Range("A140").Value = "Madagascar"
$ws.Range("A141").Value = "Liberia"
$ws.Range("A142").Value = "Etiopia"

# --- Refreshed numeric data (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos
$ws.Range("B4").Value = 990690
$ws.Range("C4").Value = 3530
$ws.Range("D4").Value = 119303
$ws.Range("E4").Value = 815534
$ws.Range("G4").Value = 440
$ws.Range("H4").Value = 55853

# Row 31 (now Chile)
$ws.Range("B31").Value = 13813
$ws.Range("C31").Value = 482
$ws.Range("D31").Value = 7327
$ws.Range("E31").Value = 6288
$ws.Range("F31").Value = 426
$ws.Range("G31").Value = 9
$ws.Range("H31").Value = 198

# Row 32 (now Japon)
$ws.Range("B32").Value = 13441
$ws.Range("D32").Value = 1809
$ws.Range("E32").Value = 11260
$ws.Range("F32").Value = 296
$ws.Range("H32").Value = 372

# Polonia
$ws.Range("B33").Value = 11902
$ws.Range("C33").Value = 285
$ws.Range("E33").Value = 8874
$ws.Range("G33").Value = 27
$ws.Range("H33").Value = 562

# Chequia
$ws.Range("B45").Value = 7431
$ws.Range("C45").Value = 27
$ws.Range("D45").Value = 2826
$ws.Range("E45").Value = 4383
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 222

# Barein
$ws.Range("B63").Value = 2708
$ws.Range("C63").Value = 61
$ws.Range("D63").Value = 1218
$ws.Range("E63").Value = 1482

# Grecia
$ws.Range("B65").Value = 2534
$ws.Range("C65").Value = 17
$ws.Range("E65").Value = 1821
$ws.Range("F65").Value = 43
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 136

# Uzbekistan
$ws.Range("D68").Value = 892
$ws.Range("E68").Value = 987

# Montenegro
$ws.Range("D123").Value = 189
$ws.Range("E123").Value = 125

# Paraguay
$ws.Range("D129").Value = 93
$ws.Range("E129").Value = 126

# Row 140 (now Madagascar)
$ws.Range("B140").Value = 128
$ws.Range("C140").Value = 4
$ws.Range("D140").Value = 75
$ws.Range("E140").Value = 53
$ws.Range("F140").Value = 1
$ws.Range("H140").Value = 0

# Row 141 (now Liberia)
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 25
$ws.Range("E141").Value = 87
$ws.Range("H141").Value = 12

# Row 142 (now Etiopia)
$ws.Range("C142").Value = 1
$ws.Range("D142").Value = 50
$ws.Range("E142").Value = 71
$ws.Range("F142").Value = 0
$ws.Range("H142").Value = 3
